$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        A = "Creary & Locke_2022_OrgSci_Breaking hte Cycle of Overwork and Recuperation.pdf"
        B = "F1_P13_Creary & Locke_2022_OrgSci_Breaking hte Cycle of Overwork and Recuperation.png"
        C = 1
        D = 2022
        E = 14
        F = 0
        G = $false
    },
    @{
        A = "Anthony_2021_ASQ_When Knowledge Work and Analytical Technology Collide.pdf"
        B = "F2_P7_Anthony_2021_ASQ_When Knowledge Work and Analytical Technology Collide.png"
        C = 2
        D = 2021
        E = 8
        F = 0
        G = $false
    },
    @{
        A = "Aoki_2020_AMJ_The Roles of Material Artfifacts in Managing the Learning Performance Paradox.pdf"
        B = "F1_P12_Aoki_2020_AMJ_The Roles of Material Artfifacts in Managing the Learning Performance Paradox.png"
        C = 1
        D = 2020
        E = 13
        F = 0
        G = $false
    },
    @{
        A = "Anthony_2021_ASQ_When Knowledge Work and Analytical Technology Collide.pdf"
        B = "F2_P21_Anthony_2021_ASQ_When Knowledge Work and Analytical Technology Collide.png"
        C = 2
        D = 2021
        E = 22
        F = 0
        G = $false
    }
)

$startRow = 145
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
}
